$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.026.43"
$ws.Range("E2").Value = "  -1.68%  "
$ws.Range("D3").Value = "1.977.07"
$ws.Range("E3").Value = "  -4.01%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.598"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.45%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.29"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.56%  "
$ws.Range("E9").Value = "  -5.54%  "
$ws.Range("E10").Value = "  -7.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0990"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.43%  "
$ws.Range("D12").Value = "2.270.26"
$ws.Range("E12").Value = "  -3.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.754"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -10.01%  "
$ws.Range("E16").Value = "  -7.19%  "
$ws.Range("D17").Value = "1.994.34"
$ws.Range("E17").Value = "  -3.19%  "
$ws.Range("D18").Value = "36.830.85"
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("E20").Value = "  -6.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "227.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("E22").Value = "  -7.01%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -9.68%  "
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("E27").Value = "  -7.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.127"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.85%  "
$ws.Range("E30").Value = "  -7.52%  "
$ws.Range("E31").Value = "  -3.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.58%  "
$ws.Range("E33").Value = "  -9.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.65%  "
$ws.Range("E35").Value = "  -7.07%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.08%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "1.418.20"
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.74%  "
$ws.Range("E43").Value = "  -8.46%  "
$ws.Range("E44").Value = "  -7.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.995"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.78%  "
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -11.86%  "
$ws.Range("D50").Value = "2.162.22"
$ws.Range("E50").Value = "  -3.51%  "
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.80%  "
